$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# employee_id, employee_name, department, absence_reason, absence_duration, absence_date, salary
$data = @{
    2  = @(59392, "João Cardoso", "Engenharia", "Doenca", 3, 45103, 7739.63)
    3  = @(76141, "Laura Rocha", "Operacoes", "Outros", 4, 45084, 9121.440000000001)
    4  = @(57370, "Nina Correia", "Financeiro", "Consulta medica", 1, 45080, 8991.209999999999)
    5  = @(24787, "Cauã Oliveira", "Atendimento ao Cliente", "Viagem de negocios", 5, 45083, 7812.4)
    6  = @(58424, "Dra. Laura Albuquerque", "Operacoes", "Consulta medica", 6, 45089, 7949.92)
    7  = @(99744, "Helena Aparecida", "Recursos Humanos", "Doenca", 3, 45097, 2861.06)
    8  = @(96217, "Rael Mendonça", "P&D", "Problemas pessoais", 2, 45094, 3876.66)
    9  = @(43150, "Arthur Gabriel Guerra", "Marketing", "Doenca", 1, 45082, 2892.67)
    10 = @(3783, "Yago da Cruz", "Marketing", "Outros", 4, 45084, 4271.39)
    11 = @(58145, "Ana Julia Jesus", "Marketing", "Viagem de negocios", 1, 45080, 8895.389999999999)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
}
